$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add D39 value
$ws.Range("D39").Value = -3.8

# Add new row 42: EndophilinA1 / 1 / SVAGLKKQFHKATQKVSEKV / -3.444
# (Set C42 before A42 so shared strings are appended in the same order as the target file)
$ws.Range("C42").Value = "SVAGLKKQFHKATQKVSEKV"
$ws.Range("A42").Value = "EndophilinA1"
$ws.Range("B42").Value = 1
$ws.Range("D42").Value = -3.444

# Update view: top left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("E42").Select()
